$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q3" sheet right after "总计" (so it becomes sheet #2,
#    pushing the existing quarter sheets down by one position each).
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $zongji)
$newSheet.Name = "2022-Q3"

# Reference sheet to copy header/row styling from (the old "2022-Q2" sheet,
# which after the insert above is now at position 3).
$styleSrc = $wb.Worksheets.Item(3)

# Copy the header row (B1:H1) formatting onto the new sheet, then fill in
# the header text.
$styleSrc.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the "index" column (A) style from the reference sheet so the new
# sheet's A column matches (bold/centered/bordered like the other sheets).
$styleSrc.Range("A2").Copy()
$newSheet.Range("A2:A16").PasteSpecial(-4122)

# Every other quarter sheet in this workbook stores B..G (fund code, name,
# scale, position, ratio, value) as plain text, even when the value looks
# numeric (e.g. "008515", "450009", "35.14") -- format the range as Text so
# entering those strings does not get auto-coerced into numbers/leading
# zeros are not dropped.
$newSheet.Range("B2:G16").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Fund holding rows for 2022-Q3 (A=index, B=code, C=name, D=scale,
# E=stock position, F=position ratio, G=held value, H=position rank).
# All of B..G are written as text (matching the convention used by every
# other quarter sheet in this workbook), A/H are written as numbers.
# ---------------------------------------------------------------------------
$rows = @(
  @(0, "450009", "国富中小盘股票", "35.14", "85.39", "3.62", "1.2721", 6),
  @(1, "450002", "国富弹性市值混合", "30.14", "85.83", "3.18", "0.9585", 8),
  @(2, "012239", "惠升优势企业一年持有期灵活配置混合", "10.70", "79.97", "3.15", "0.3370", 7),
  @(3, "008515", "富兰克林国海基本面优选混合", "11.82", "85.30", "2.33", "0.2754", 10),
  @(4, "159916", "建信深证基本面60ETF", "3.55", "98.81", "2.62", "0.0930", 10),
  @(5, "159910", "嘉实深证基本面120ETF", "3.42", "99.58", "2.10", "0.0718", 10),
  @(6, "000058", "国联安安泰灵活配置混合", "4.16", "31.99", "1.11", "0.0462", 9),
  @(7, "673020", "西部利得成长精选灵活配置混合", "1.42", "86.11", "2.85", "0.0405", 10),
  @(8, "004131", "国联安鑫发混合A", "3.26", "23.06", "1.09", "0.0355", 7),
  @(9, "002186", "国联安鑫享灵活配置混合C", "1.64", "33.73", "1.25", "0.0205", 7),
  @(10, "159913", "交银深证300价值ETF", "0.40", "97.69", "2.14", "0.0086", 9),
  @(11, "001228", "国联安鑫享灵活配置混合A", "0.67", "33.73", "1.25", "0.0084", 7),
  @(12, "004132", "国联安鑫发混合C", "0.57", "23.06", "1.09", "0.0062", 7),
  @(13, "007288", "合煦智远消费主题股票C", "0.47", "20.37", "0.93", "0.0044", 8),
  @(14, "007287", "合煦智远消费主题股票A", "0.05", "20.37", "0.93", "0.0005", 8)
)

$r = 2
foreach ($row in $rows) {
  $newSheet.Cells.Item($r, 1).Value = $row[0]
  $newSheet.Cells.Item($r, 2).Value = $row[1]
  $newSheet.Cells.Item($r, 3).Value = $row[2]
  $newSheet.Cells.Item($r, 4).Value = $row[3]
  $newSheet.Cells.Item($r, 5).Value = $row[4]
  $newSheet.Cells.Item($r, 6).Value = $row[5]
  $newSheet.Cells.Item($r, 7).Value = $row[6]
  $newSheet.Cells.Item($r, 8).Value = $row[7]
  $r = $r + 1
}

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row 2 for 2022-Q3 and bump
#    the running index in column A for every pre-existing row.
# ---------------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()
$zongji.Rows.Item(2).ClearFormats()

$zongji.Cells.Item(3, 1).Copy()
$zongji.Cells.Item(2, 1).PasteSpecial(-4122)

$zongji.Cells.Item(2, 1).Value = 0
$zongji.Cells.Item(2, 2).Value = "2022-Q3"
$zongji.Cells.Item(2, 3).Value = 15
$zongji.Cells.Item(2, 4).Value = 3.18

for ($row = 3; $row -le 9; $row++) {
  $zongji.Cells.Item($row, 1).Value = $row - 2
}

$zongji.Range("A1").Select()
